{"js": "// Locate the paragraph containing \"Tim Gleeson 17664081\" (it currently\n// carries the hidden \"_GoBack\" bookmark at its end).\nconst timResults = context.document.body.search(\"Tim Gleeson 17664081\", { matchCase: true });\ntimResults.load(\"items\");\nawait context.sync();\n\nconst timPara = timResults.items[0].paragraphs.getFirst();\n\n// Insert the two new team-member paragraphs right after Tim Gleeson's line.\nconst fuadPara = timPara.insertParagraph(\"Fuad Faraj S Aljohani 17619583\", Word.InsertLocation.after);\nawait context.sync();\n\nfuadPara.insertParagraph(\"Hassan Iqbal 18026141\", Word.InsertLocation.after);\nawait context.sync();\n\n// The \"_GoBack\" bookmark needs to move from the end of Tim Gleeson's\n// paragraph to the end of the new last paragraph (\"Hassan Iqbal ...\").\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// Re-locate the freshly-inserted paragraph by search: reusing the object\n// returned from insertParagraph() for a follow-up getRange() call can\n// target a stale position, so look it up again after the sync above.\nconst hassanResults = context.document.body.search(\"Hassan Iqbal 18026141\", { matchCase: true });\nhassanResults.load(\"items\");\nawait context.sync();\n\nconst hassanPara = hassanResults.items[0].paragraphs.getFirst();\nconst hassanEnd = hassanPara.getRange(Word.RangeLocation.end);\nhassanEnd.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Add two missing team members (\"Fuad Faraj S Aljohani\" and \"Hassan Iqbal\")\n# to the member list, right after \"Tim Gleeson 17664081\". The hidden\n# \"_GoBack\" bookmark that currently sits at the end of Tim Gleeson's line\n# must end up at the end of the new last line (\"Hassan Iqbal 18026141\").\n\n$d = $word.ActiveDocument\n\n# --- Locate \"Tim Gleeson 17664081\" -----------------------------------\n$timRange = $d.Content\n$timRange.Find.Execute(\"Tim Gleeson 17664081\") | Out-Null\n$timPara = $timRange.Paragraphs(1)\n\n# --- Remove the old \"_GoBack\" bookmark up front -----------------------\n$goBack = $d.Bookmarks(\"_GoBack\")\n$goBack.Delete()\n\n# --- Insert the two new paragraphs -------------------------------------\n$timPara.Range.InsertParagraphAfter()\n$fuadPara = $d.Paragraphs($timPara.Index + 1)\n$fuadPara.Range.InsertAfter(\"Fuad Faraj S Aljohani 17619583\")\n\n$fuadPara = $d.Paragraphs($timPara.Index + 1)\n$fuadPara.Range.InsertParagraphAfter()\n$hassanPara = $d.Paragraphs($timPara.Index + 2)\n$hassanPara.Range.InsertAfter(\"Hassan Iqbal 18026141\")\n\n# --- Re-create \"_GoBack\" at the end of the \"Hassan Iqbal\" paragraph ----\n# NOTE: adding a bookmark with a zero-width (collapsed) range placed\n# exactly one character before a paragraph mark lands it in the wrong\n# spot in this host, so a one-character placeholder is appended first,\n# the bookmark is anchored just before it (a position that is *not* a\n# paragraph-end boundary), and the placeholder is then deleted. The\n# bookmark stays correctly anchored after the deletion.\n$hassanPara = $d.Paragraphs($timPara.Index + 2)\n$hassanPara.Range.InsertAfter(\"X\")\n$hassanPara = $d.Paragraphs($timPara.Index + 2)\n$bookmarkPos = $hassanPara.Range.End - 2\n$d.Bookmarks.Add(\"_GoBack\", $d.Range($bookmarkPos, $bookmarkPos))\n$d.Range($bookmarkPos, $bookmarkPos + 1).Delete()\n"}
